$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for area1 vs rest-of-areas cluster split
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 257

# Remove the now-unused third cluster row (row 4: A4=2, B4=94)
$ws.Range("A4:B4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
